$d = $word.ActiveDocument
$s = $d.Shapes.Item(1)

$s.Left   = -0.499921259842520
$s.Top    = 31.3507874015748
$s.Width  = 648.0
$s.Height = 334.973307086614
